$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    7.273612799567631,
    -7.71094665953368,
    -2.953243576515993,
    2.092325086009379,
    -0.8957659543836369,
    0.808111282324586,
    -0.5400899636021528,
    0.9306312308838732,
    -0.0521044692539093,
    -1.185636841993183,
    0.7007593354088806,
    0.876933787738463,
    0.3962997177008458,
    -0.062173517,
    0.2546725873431263,
    0.3607736698111976,
    0.08648729331993882,
    -0.5990062151031799,
    0.6838772855378736,
    -0.2042791748566441,
    0.5947585845621888
)

$startRow = 3
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
